# Applies the "secdep collection sheet scenarios" edit:
#  - Adds a ReceiptNumber row (with receipt number value) above the existing
#    "submitmakerepayment" row on the "Modify Transaction1/2/3" sheets.
#  - Updates the selection / scroll position on a few sheets.

$wb = $excel.ActiveWorkbook

# --- Modify Transaction1 (receipt number 7654) ---
$ws = $wb.Worksheets.Item("Modify Transaction1")
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "ReceiptNumber"
$ws.Range("B4").Value = 7654
$ws.Range("B8").Select()

# --- Modify Transaction2 (receipt number 32) ---
$ws = $wb.Worksheets.Item("Modify Transaction2")
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "ReceiptNumber"
$ws.Range("B4").Value = 32
$ws.Range("B10").Select()

# --- Modify Transaction3 (receipt number 76543) ---
$ws = $wb.Worksheets.Item("Modify Transaction3")
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "ReceiptNumber"
$ws.Range("B4").Value = 76543
$ws.Range("B9").Select()

# --- Transactions sheet: scroll so row 7 is at the top ---
$ws = $wb.Worksheets.Item("Transactions")
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
